# NIT-9004967659.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# - Removes all rows for TANIA MARCELA ESPINOSA TENORIO (1143366477), periods 2506-2508
# - Adds a brand-new period 2509 for the 4 recurring workers, replacing TANIA with
#   BAYRON EMIRO CONSUEGRA OSUNA (1148449763)
# - Updates the "VALOR MORA" and "Cant. Periodos" summary cells accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Delete TANIA MARCELA ESPINOSA TENORIO's three rows (periods 2506, 2507,
#    2508). Deleting bottom-up keeps the earlier row numbers stable.
# ---------------------------------------------------------------------------
$ws.Rows("50:50").Delete()
$ws.Rows("46:46").Delete()
$ws.Rows("42:42").Delete()

# ---------------------------------------------------------------------------
# 2) Insert four new rows (period 2509) after the current last data row (49),
#    copying row 49 as a formatting template one row at a time so borders /
#    number formats carry through correctly, then reapply the template format
#    to be safe.
# ---------------------------------------------------------------------------
$ws.Rows("49:49").Copy()
$ws.Rows("50:50").Insert()
$ws.Rows("50:50").Copy()
$ws.Rows("51:51").Insert()
$ws.Rows("51:51").Copy()
$ws.Rows("52:52").Insert()
$ws.Rows("52:52").Copy()
$ws.Rows("53:53").Insert()

$ws.Range("B49:J49").Copy()
$ws.Range("B50:J53").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Fill in the new period-2509 data (same pattern used by 2506-2508):
#    HEIBER TOBIO PEREZ, BAYRON EMIRO CONSUEGRA OSUNA (new), ENA ONIX SANTOS
#    GUTIERREZ, DAVID ENRIQUE OSPINO CONRADO.
# ---------------------------------------------------------------------------
$ws.Range("B50").Value = "CC"
$ws.Range("C50").Value = "1047386799"
$ws.Range("D50").Value = "HEIBER TOBIO PEREZ"
$ws.Range("E50").Value = "2509"
$ws.Range("F50").Value = 56940
$ws.Range("G50").Value = 1423500

$ws.Range("B51").Value = "CC"
$ws.Range("C51").Value = "1148449763"
$ws.Range("D51").Value = "BAYRON EMIRO CONSUEGRA OSUNA"
$ws.Range("E51").Value = "2509"
$ws.Range("F51").Value = 56940
$ws.Range("G51").Value = 1423500

$ws.Range("B52").Value = "CC"
$ws.Range("C52").Value = "1072248349"
$ws.Range("D52").Value = "ENA ONIX SANTOS GUTIERREZ"
$ws.Range("E52").Value = "2509"
$ws.Range("F52").Value = 26400
$ws.Range("G52").Value = 660000

$ws.Range("B53").Value = "CC"
$ws.Range("C53").Value = "1065829293"
$ws.Range("D53").Value = "DAVID ENRIQUE OSPINO CONRADO"
$ws.Range("E53").Value = "2509"
$ws.Range("F53").Value = 56940
$ws.Range("G53").Value = 1423500

# Center-align the "Periodo Mora" column like the rest of the refreshed table
$ws.Range("E16:E53").HorizontalAlignment = -4108  # xlCenter

# ---------------------------------------------------------------------------
# 4) Update the summary header: total "VALOR MORA" grows by the new row's
#    26400, and "Cant. Periodos" grows from 20 to 21. "Cant. Trabajadores"
#    stays at 6 (TANIA swapped 1-for-1 with BAYRON).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1443259
$ws.Range("F13").Value = 21
